$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text representation (with trailing
# zeros / full precision) instead of Excel auto-converting numeric-looking
# strings ("1.001", "322.90", ...) into floating point numbers.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D22","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '28.047.17'
$ws.Range("E2").Value = '  +0.52%  '

# Row 3
$ws.Range("D3").Value = '1.765.56'
$ws.Range("E3").Value = '  -2.08%  '

# Row 4
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.71%  '

# Row 5
$ws.Range("D5").Value = '322.90'
$ws.Range("E5").Value = '  -0.55%  '

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.96%  '

# Row 7
$ws.Range("D7").Value = '0.4243'
$ws.Range("E7").Value = '  -4.49%  '

# Row 8
$ws.Range("D8").Value = '0.3603'
$ws.Range("E8").Value = '  -3.99%  '

# Row 9
$ws.Range("D9").Value = '44.33'
$ws.Range("E9").Value = '  -1.44%  '

# Row 10
$ws.Range("D10").Value = '0.07452'
$ws.Range("E10").Value = '  -3.55%  '

# Row 11
$ws.Range("D11").Value = '1.102'
$ws.Range("E11").Value = '  -2.64%  '

# Row 12
$ws.Range("D12").Value = '0.9993'
$ws.Range("E12").Value = '  +0.89%  '

# Row 13
$ws.Range("D13").Value = '21.38'
$ws.Range("E13").Value = '  -2.64%  '

# Row 14
$ws.Range("D14").Value = '6.090'
$ws.Range("E14").Value = '  -2.27%  '

# Row 15
$ws.Range("D15").Value = '7.350'
$ws.Range("E15").Value = '  -1.73%  '

# Row 16
$ws.Range("D16").Value = '1.795.56'
$ws.Range("E16").Value = '  +1.14%  '

# Row 17
$ws.Range("D17").Value = '91.26'
$ws.Range("E17").Value = '  -0.17%  '

# Row 18
$ws.Range("D18").Value = '0.00001061'
$ws.Range("E18").Value = '  -1.52%  '

# Row 19
$ws.Range("D19").Value = '0.06418'
$ws.Range("E19").Value = '  +1.38%  '

# Row 20
$ws.Range("D20").Value = '0.9988'
$ws.Range("E20").Value = '  +0.87%  '

# Row 21
$ws.Range("D21").Value = '17.08'
$ws.Range("E21").Value = '  -2.58%  '

# Row 22
$ws.Range("D22").Value = '5.976'
$ws.Range("E22").Value = '  -4.85%  '

# Row 23
$ws.Range("D23").Value = '28.026.35'
$ws.Range("E23").Value = '  +0.28%  '

# Row 24
$ws.Range("D24").Value = '11.29'
$ws.Range("E24").Value = '  -3.02%  '

# Row 25
$ws.Range("D25").Value = '2.148'
$ws.Range("E25").Value = '  -2.14%  '

# Row 26
$ws.Range("D26").Value = '158.25'
$ws.Range("E26").Value = '  +3.58%  '

# Row 27
$ws.Range("D27").Value = '20.16'
$ws.Range("E27").Value = '  -3.32%  '

# Row 28
$ws.Range("D28").Value = '1.992.01'
$ws.Range("E28").Value = '  +1.11%  '

# Row 29
$ws.Range("D29").Value = '2.138'
$ws.Range("E29").Value = '  -8.90%  '

# Row 30
$ws.Range("D30").Value = '126.38'
$ws.Range("E30").Value = '  -1.53%  '

# Row 31
$ws.Range("D31").Value = '1.173'
$ws.Range("E31").Value = '  -2.33%  '

# Row 32
$ws.Range("D32").Value = '5.643'
$ws.Range("E32").Value = '  -2.64%  '

# Row 33
$ws.Range("D33").Value = '0.09037'
$ws.Range("E33").Value = '  -1.96%  '

# Row 34
$ws.Range("D34").Value = '3.527'
$ws.Range("E34").Value = '  -3.06%  '

# Row 35
$ws.Range("D35").Value = '12.56'
$ws.Range("E35").Value = '  -1.24%  '

# Row 36
$ws.Range("D36").Value = '0.02327'
$ws.Range("E36").Value = '  -0.21%  '

# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2106'
$ws.Range("E37").Value = '  -3.34%  '

# Row 38
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '5.038'
$ws.Range("E38").Value = '  -1.96%  '

# Row 39
$ws.Range("D39").Value = '0.06052'
$ws.Range("E39").Value = '  -1.81%  '

# Row 40
$ws.Range("D40").Value = '0.6401'
$ws.Range("E40").Value = '  -2.26%  '

# Row 41
$ws.Range("D41").Value = '1.187'
$ws.Range("E41").Value = '  +0.35%  '

# Row 42
$ws.Range("D42").Value = '0.9994'
$ws.Range("E42").Value = '  +0.69%  '

# Row 43
$ws.Range("D43").Value = '7.821'
$ws.Range("E43").Value = '  -2.84%  '

# Row 44
$ws.Range("D44").Value = '1.391'
$ws.Range("E44").Value = '  +0.71%  '

# Row 45
$ws.Range("D45").Value = '13.54'
$ws.Range("E45").Value = '  -1.88%  '

# Row 46
$ws.Range("D46").Value = '0.5947'
$ws.Range("E46").Value = '  -1.88%  '

# Row 47
$ws.Range("E47").Value = '  -0.67%  '

# Row 48
$ws.Range("D48").Value = '2.022'
$ws.Range("E48").Value = '  +0.68%  '

# Row 49
$ws.Range("D49").Value = '123.14'
$ws.Range("E49").Value = '  -2.76%  '

# Row 50
$ws.Range("D50").Value = '1.193'
$ws.Range("E50").Value = '  +4.51%  '

# Row 51
$ws.Range("E51").Value = '  -0.90%  '
